$wb = $excel.ActiveWorkbook
$wsMonitor = $wb.Worksheets.Item("Monitor")
$wsHoldings = $wb.Worksheets.Item("Current_Holdings")

# --- Monitor sheet: rows 5-9 get rotated (row7,row8,row9,row5,row6 -> rows 5..9) ---
# Capture the current (pre-edit) row values before overwriting anything.
$origB5 = $wsMonitor.Range("B5").Value2
$origC5 = $wsMonitor.Range("C5").Value2
$origD5 = $wsMonitor.Range("D5").Value2
$origE5 = $wsMonitor.Range("E5").Value2
$origF5 = $wsMonitor.Range("F5").Value2
$origG5 = $wsMonitor.Range("G5").Value2
$origI5 = $wsMonitor.Range("I5").Value2
$origK5 = $wsMonitor.Range("K5").Value2
$origL5 = $wsMonitor.Range("L5").Value2
$origM5 = $wsMonitor.Range("M5").Value2

$origB6 = $wsMonitor.Range("B6").Value2
$origC6 = $wsMonitor.Range("C6").Value2
$origD6 = $wsMonitor.Range("D6").Value2
$origE6 = $wsMonitor.Range("E6").Value2
$origF6 = $wsMonitor.Range("F6").Value2
$origG6 = $wsMonitor.Range("G6").Value2
$origI6 = $wsMonitor.Range("I6").Value2
$origK6 = $wsMonitor.Range("K6").Value2
$origL6 = $wsMonitor.Range("L6").Value2
$origM6 = $wsMonitor.Range("M6").Value2

$origB7 = $wsMonitor.Range("B7").Value2
$origC7 = $wsMonitor.Range("C7").Value2
$origD7 = $wsMonitor.Range("D7").Value2
$origE7 = $wsMonitor.Range("E7").Value2
$origF7 = $wsMonitor.Range("F7").Value2
$origG7 = $wsMonitor.Range("G7").Value2
$origI7 = $wsMonitor.Range("I7").Value2
$origK7 = $wsMonitor.Range("K7").Value2
$origL7 = $wsMonitor.Range("L7").Value2
$origM7 = $wsMonitor.Range("M7").Value2

$origB8 = $wsMonitor.Range("B8").Value2
$origC8 = $wsMonitor.Range("C8").Value2
$origD8 = $wsMonitor.Range("D8").Value2
$origE8 = $wsMonitor.Range("E8").Value2
$origF8 = $wsMonitor.Range("F8").Value2
$origG8 = $wsMonitor.Range("G8").Value2
$origI8 = $wsMonitor.Range("I8").Value2
$origK8 = $wsMonitor.Range("K8").Value2
$origL8 = $wsMonitor.Range("L8").Value2
$origM8 = $wsMonitor.Range("M8").Value2

$origB9 = $wsMonitor.Range("B9").Value2
$origC9 = $wsMonitor.Range("C9").Value2
$origD9 = $wsMonitor.Range("D9").Value2
$origE9 = $wsMonitor.Range("E9").Value2
$origF9 = $wsMonitor.Range("F9").Value2
$origG9 = $wsMonitor.Range("G9").Value2
$origI9 = $wsMonitor.Range("I9").Value2
$origK9 = $wsMonitor.Range("K9").Value2
$origL9 = $wsMonitor.Range("L9").Value2
$origM9 = $wsMonitor.Range("M9").Value2

# New row 5 <- old row 7, with a couple of figures corrected (the "figures_in" bug fix)
$wsMonitor.Range("B5").Value = $origB7
$wsMonitor.Range("C5").Value = $origC7
$wsMonitor.Range("D5").Value = $origD7
$wsMonitor.Range("E5").Value = 1.73
$wsMonitor.Range("F5").Value = 0.11290796645794088
$wsMonitor.Range("G5").Value = 0.09786536075997788
$wsMonitor.Range("I5").Value = $origI7
$wsMonitor.Range("K5").Value = $origK7
$wsMonitor.Range("L5").Value = $origL7
$wsMonitor.Range("M5").Value = $origM7

# New row 6 <- old row 8 (unchanged figures)
$wsMonitor.Range("B6").Value = $origB8
$wsMonitor.Range("C6").Value = $origC8
$wsMonitor.Range("D6").Value = $origD8
$wsMonitor.Range("E6").Value = $origE8
$wsMonitor.Range("F6").Value = $origF8
$wsMonitor.Range("G6").Value = $origG8
$wsMonitor.Range("I6").Value = $origI8
$wsMonitor.Range("K6").Value = $origK8
$wsMonitor.Range("L6").Value = $origL8
$wsMonitor.Range("M6").Value = $origM8

# New row 7 <- old row 9 (unchanged figures)
$wsMonitor.Range("B7").Value = $origB9
$wsMonitor.Range("C7").Value = $origC9
$wsMonitor.Range("D7").Value = $origD9
$wsMonitor.Range("E7").Value = $origE9
$wsMonitor.Range("F7").Value = $origF9
$wsMonitor.Range("G7").Value = $origG9
$wsMonitor.Range("I7").Value = $origI9
$wsMonitor.Range("K7").Value = $origK9
$wsMonitor.Range("L7").Value = $origL9
$wsMonitor.Range("M7").Value = $origM9

# New row 8 <- old row 5 (unchanged figures)
$wsMonitor.Range("B8").Value = $origB5
$wsMonitor.Range("C8").Value = $origC5
$wsMonitor.Range("D8").Value = $origD5
$wsMonitor.Range("E8").Value = $origE5
$wsMonitor.Range("F8").Value = $origF5
$wsMonitor.Range("G8").Value = $origG5
$wsMonitor.Range("I8").Value = $origI5
$wsMonitor.Range("K8").Value = $origK5
$wsMonitor.Range("L8").Value = $origL5
$wsMonitor.Range("M8").Value = $origM5

# New row 9 <- old row 6, with the "last update" date corrected
$wsMonitor.Range("B9").Value = $origB6
$wsMonitor.Range("C9").Value = $origC6
$wsMonitor.Range("D9").Value = $origD6
$wsMonitor.Range("E9").Value = $origE6
$wsMonitor.Range("F9").Value = $origF6
$wsMonitor.Range("G9").Value = $origG6
$wsMonitor.Range("I9").Value = $origI6
$wsMonitor.Range("K9").ClearContents()
$wsMonitor.Range("L9").Value = 44903
$wsMonitor.Range("M9").Value = $origM6

# --- Current_Holdings sheet: last-update date correction ---
$wsHoldings.Range("I2").Value = 44904

# --- Selection / active-sheet bookkeeping ---
$wsMonitor.Range("B5:R17").Select()
$wsHoldings.Activate()
$wsHoldings.Range("B7:K13").Select()
